$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 19.05064621328966
$ws.Range("D2").Value = 4.606291277960393
$ws.Range("E2").Value = 10.4977862042877
$ws.Range("F2").Value = 69.07147766769589
$ws.Range("G2").Value = 3.833506227394755
$ws.Range("J2").Value = 10.96933330052918
$ws.Range("K2").Value = 16.93207738241814
$ws.Range("L2").Value = 11.23066201579071
$ws.Range("M2").Value = 18.17708175807337
$ws.Range("B3").Value = 19.07214369751059
$ws.Range("D3").Value = 4.444431487095964
$ws.Range("E3").Value = 10.47249905381721
$ws.Range("F3").Value = 68.25161133668989
$ws.Range("G3").Value = 3.837460910741264
$ws.Range("J3").Value = 10.94651261277166
$ws.Range("K3").Value = 16.94663875645211
$ws.Range("L3").Value = 11.28428004849743
$ws.Range("M3").Value = 18.23747938262637
$ws.Range("B4").Value = 19.09012157755582
$ws.Range("D4").Value = 4.341587258140049
$ws.Range("E4").Value = 10.45645879720313
$ws.Range("F4").Value = 67.74649462972539
$ws.Range("G4").Value = 3.840013245919101
$ws.Range("J4").Value = 10.93201593665261
$ws.Range("K4").Value = 16.96189442160391
$ws.Range("L4").Value = 11.31978788780053
$ws.Range("M4").Value = 18.27822904561163
$ws.Range("B5").Value = 19.09864844442745
$ws.Range("D5").Value = 4.298847363366511
$ws.Range("E5").Value = 10.44979239286988
$ws.Range("F5").Value = 67.54036370310389
$ws.Range("G5").Value = 3.841084683237937
$ws.Range("J5").Value = 10.92598589380079
$ws.Range("K5").Value = 16.96969521198523
$ws.Range("L5").Value = 11.33490761401541
$ws.Range("M5").Value = 18.29575609060415
$ws.Range("B6").Value = 19.10013681221506
$ws.Range("D6").Value = 4.291701570383426
$ws.Range("E6").Value = 10.44867757842604
$ws.Range("F6").Value = 67.50612226302233
$ws.Range("G6").Value = 3.84126449106095
$ws.Range("J6").Value = 10.92497719015007
$ws.Range("K6").Value = 16.97108605754444
$ws.Range("L6").Value = 11.33745748224175
$ws.Range("M6").Value = 18.29872206948848
$ws.Range("B7").Value = 19.09023171362632
$ws.Range("D7").Value = 4.341014157120925
$ws.Range("E7").Value = 10.45636941841145
$ws.Range("F7").Value = 67.74371567366774
$ws.Range("G7").Value = 3.84002756864197
$ws.Range("J7").Value = 10.93193511033626
$ws.Range("K7").Value = 16.96199321848963
$ws.Range("L7").Value = 11.31998916636263
$ws.Range("M7").Value = 18.27846169211156
$ws.Range("B8").Value = 19.05706645022899
$ws.Range("D8").Value = 4.551223420148466
$ws.Range("E8").Value = 10.48917344277654
$ws.Range("F8").Value = 68.78922765604241
$ws.Range("G8").Value = 3.834844110472071
$ws.Range("J8").Value = 10.96156492814147
$ws.Range("K8").Value = 16.93578545318809
$ws.Range("L8").Value = 11.24861256400696
$ws.Range("M8").Value = 18.19714587725937
$ws.Range("B9").Value = 19.02996116058729
$ws.Range("D9").Value = 4.934415278699615
$ws.Range("E9").Value = 10.54945599882478
$ws.Range("F9").Value = 70.81974224764122
$ws.Range("G9").Value = 3.825658769891833
$ws.Range("J9").Value = 11.01585358876082
$ws.Range("K9").Value = 16.93462148849137
$ws.Range("L9").Value = 11.12917817935402
$ws.Range("M9").Value = 18.0667872607765
$ws.Range("B10").Value = 19.03316746152063
$ws.Range("D10").Value = 5.196441725259515
$ws.Range("E10").Value = 10.59132400217888
$ws.Range("F10").Value = 72.2914055903917
$ws.Range("G10").Value = 3.819499502029549
$ws.Range("J10").Value = 11.05345704403104
$ws.Range("K10").Value = 16.96447687087403
$ws.Range("L10").Value = 11.05397020349245
$ws.Range("M10").Value = 17.98877921055272
$ws.Range("B11").Value = 19.03963492356994
$ws.Range("D11").Value = 5.311100072325751
$ws.Range("E11").Value = 10.60985477858281
$ws.Range("F11").Value = 72.95483081480826
$ws.Range("G11").Value = 3.816823747751932
$ws.Range("J11").Value = 11.07007758504955
$ws.Range("K11").Value = 16.9847197176867
$ws.Range("L11").Value = 11.02248446787241
$ws.Range("M11").Value = 17.95715397526462
$ws.Range("B12").Value = 19.04280230026004
$ws.Range("D12").Value = 5.353843994926468
$ws.Range("E12").Value = 10.61679855032455
$ws.Range("F12").Value = 73.20504416906512
$ws.Range("G12").Value = 3.81582851660652
$ws.Range("J12").Value = 11.07630224647482
$ws.Range("K12").Value = 16.99334023732566
$ws.Range("L12").Value = 11.01095426879745
$ws.Range("M12").Value = 17.94573388176347
$ws.Range("B13").Value = 19.04208824052167
$ws.Range("D13").Value = 5.344668653954096
$ws.Range("E13").Value = 10.61530634031536
$ws.Range("F13").Value = 73.15120346262842
$ws.Range("G13").Value = 3.816042057919828
$ws.Range("J13").Value = 11.07496472088683
$ws.Range("K13").Value = 16.99144124132398
$ws.Range("L13").Value = 11.01342002428657
$ws.Range("M13").Value = 17.94816868121134
$ws.Range("B14").Value = 19.03988112335233
$ws.Range("D14").Value = 5.314630278466329
$ws.Range("E14").Value = 10.61042751478696
$ws.Range("F14").Value = 72.97543709257052
$ws.Range("G14").Value = 3.816741509035925
$ws.Range("J14").Value = 11.07059107448207
$ws.Range("K14").Value = 16.98540981503323
$ws.Range("L14").Value = 11.02152799596603
$ws.Range("M14").Value = 17.95620329769251
$ws.Range("B15").Value = 19.03862267425018
$ws.Range("D15").Value = 5.296142435803069
$ws.Range("E15").Value = 10.60742954977336
$ws.Range("F15").Value = 72.8676392285837
$ws.Range("G15").Value = 3.817172286138804
$ws.Range("J15").Value = 11.06790309964139
$ws.Range("K15").Value = 16.98183965101573
$ws.Range("L15").Value = 11.02654553249518
$ws.Range("M15").Value = 17.96119711403247
$ws.Range("B16").Value = 19.03284543280574
$ws.Range("D16").Value = 5.188855249423961
$ws.Range("E16").Value = 10.59010266118757
$ws.Range("F16").Value = 72.24791535446332
$ws.Range("G16").Value = 3.819676896772298
$ws.Range("J16").Value = 11.05236113795064
$ws.Range("K16").Value = 16.96328778024448
$ws.Range("L16").Value = 11.05608279558734
$ws.Range("M16").Value = 17.9909237303239
$ws.Range("B17").Value = 19.03058294909339
$ws.Range("D17").Value = 5.121858915090774
$ws.Range("E17").Value = 10.57934168790851
$ws.Range("F17").Value = 71.86608494274014
$ws.Range("G17").Value = 3.821245614440547
$ws.Range("J17").Value = 11.04270276018957
$ws.Range("K17").Value = 16.95361132155678
$ws.Range("L17").Value = 11.0749017990841
$ws.Range("M17").Value = 18.0101493073356
$ws.Range("B18").Value = 19.02975324435481
$ws.Range("D18").Value = 5.082898301121108
$ws.Range("E18").Value = 10.57310372963152
$ws.Range("F18").Value = 71.64590585047443
$ws.Range("G18").Value = 3.822159778113165
$ws.Range("J18").Value = 11.03710178536409
$ws.Range("K18").Value = 16.94867294901126
$ws.Range("L18").Value = 11.08598260563984
$ws.Range("M18").Value = 18.02157074194023
$ws.Range("B19").Value = 19.0295533722442
$ws.Range("D19").Value = 5.069634444672513
$ws.Range("E19").Value = 10.57098330905341
$ws.Range("F19").Value = 71.57126515027667
$ws.Range("G19").Value = 3.822471342015817
$ws.Range("J19").Value = 11.03519751785359
$ws.Range("K19").Value = 16.94710870219709
$ws.Range("L19").Value = 11.08977843014576
$ws.Range("M19").Value = 18.02550024057206
$ws.Range("B20").Value = 19.03077499646815
$ws.Range("D20").Value = 5.129035057360382
$ws.Range("E20").Value = 10.58049222994136
$ws.Range("F20").Value = 71.90679047779224
$ws.Range("G20").Value = 3.821077393216907
$ws.Range("J20").Value = 11.04373563943117
$ws.Range("K20").Value = 16.95457649390072
$ws.Range("L20").Value = 11.07287192080437
$ws.Range("M20").Value = 18.00806509856544
$ws.Range("B21").Value = 19.04050993162229
$ws.Range("D21").Value = 5.3234717498925
$ws.Range("E21").Value = 10.61186253180394
$ws.Range("F21").Value = 73.02709250457187
$ws.Range("G21").Value = 3.816535575198143
$ws.Range("J21").Value = 11.07187759286334
$ws.Range("K21").Value = 16.98715550280782
$ws.Range("L21").Value = 11.01913582319097
$ws.Range("M21").Value = 17.95382825097327
$ws.Range("B22").Value = 19.05105769325058
$ws.Range("D22").Value = 5.446606457528967
$ws.Range("E22").Value = 10.63193712879807
$ws.Range("F22").Value = 73.75331584080553
$ws.Range("G22").Value = 3.81367220813743
$ws.Range("J22").Value = 11.08986706322289
$ws.Range("K22").Value = 17.01401190060858
$ws.Range("L22").Value = 10.9863058408944
$ws.Range("M22").Value = 17.92162055321686
$ws.Range("B23").Value = 19.0450459900294
$ws.Range("D23").Value = 5.381254276152606
$ws.Range("E23").Value = 10.62126185801108
$ws.Range("F23").Value = 73.36630823662247
$ws.Range("G23").Value = 3.815190874849196
$ws.Range("J23").Value = 11.08030240356616
$ws.Range("K23").Value = 16.9991703270765
$ws.Range("L23").Value = 11.00361809054363
$ws.Range("M23").Value = 17.9385138705452
$ws.Range("B24").Value = 19.03068670456132
$ws.Range("D24").Value = 5.125792103231314
$ws.Range("E24").Value = 10.57997222982327
$ws.Range("F24").Value = 71.88838954937272
$ws.Range("G24").Value = 3.821153407751675
$ws.Range("J24").Value = 11.04326882504559
$ws.Range("K24").Value = 16.95413819331439
$ws.Range("L24").Value = 11.07378881427696
$ws.Range("M24").Value = 18.00900622192192
$ws.Range("B25").Value = 19.03323127595131
$ws.Range("D25").Value = 4.834060657083127
$ws.Range("E25").Value = 10.53357647791251
$ws.Range("F25").Value = 70.27339008761291
$ws.Range("G25").Value = 3.828039610845458
$ws.Range("J25").Value = 11.00157101785063
$ws.Range("K25").Value = 16.9295446901695
$ws.Range("L25").Value = 11.15928778315216
$ws.Range("M25").Value = 18.09893443669819
